$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(40, 8).Value = 6038.25  # H40: 5562.56 -> 6038.25
$ws.Cells.Item(40, 9).Value = 5000.625  # I40: 4700.4443 -> 5000.625
$ws.Cells.Item(40, 10).Value = 6730  # J40: 6047.5 -> 6730
$ws.Cells.Item(40, 11).Value = 5000.625  # K40: 4700.4443 -> 5000.625
$ws.Cells.Item(40, 12).Value = 6730  # L40: 6047.5 -> 6730
$ws.Cells.Item(40, 13).Value = -4825.625  # M40: -4525.4443 -> -4825.625
$ws.Cells.Item(40, 14).Value = -7080  # N40: -6397.5 -> -7080
$ws.Cells.Item(55, 8).Value = 255058  # H55: 253641.5 -> 255058
$ws.Cells.Item(55, 9).Value = 0  # I55: 1000 -> 0
$ws.Cells.Item(55, 10).Value = 255058  # J55: 337855.34 -> 255058
$ws.Cells.Item(55, 11).Value = 0  # K55: 1000 -> 0
$ws.Cells.Item(55, 12).Value = 255058  # L55: 337855.34 -> 255058
$ws.Cells.Item(55, 13).ClearContents()  # M55: -786 -> (cleared)
$ws.Cells.Item(55, 14).Value = -255486  # N55: -338283.34 -> -255486
$ws.Cells.Item(62, 8).Value = 11368222  # H62: 12504245 -> 11368222
$ws.Cells.Item(62, 10).Value = 7999  # J62: 7999.5 -> 7999
$ws.Cells.Item(62, 12).Value = 7999  # L62: 7999.5 -> 7999
$ws.Cells.Item(62, 14).Value = -9247  # N62: -9247.5 -> -9247
$ws.Cells.Item(64, 8).Value = 5089.4  # H64: 5090 -> 5089.4
$ws.Cells.Item(64, 9).Value = 4331.3335  # I64: 4332.3335 -> 4331.3335
$ws.Cells.Item(64, 11).Value = 4331.3335  # K64: 4332.3335 -> 4331.3335
$ws.Cells.Item(64, 13).Value = -4083.3335  # M64: -4084.3335 -> -4083.3335
$ws.Cells.Item(65, 8).Value = 11368222  # H65: 12504245 -> 11368222
$ws.Cells.Item(65, 10).Value = 7999  # J65: 7999.5 -> 7999
$ws.Cells.Item(65, 12).Value = 39995  # L65: 39997.5 -> 39995
$ws.Cells.Item(65, 14).Value = -46235  # N65: -46237.5 -> -46235
$ws.Cells.Item(67, 8).Value = 5089.4  # H67: 5090 -> 5089.4
$ws.Cells.Item(67, 9).Value = 4331.3335  # I67: 4332.3335 -> 4331.3335
$ws.Cells.Item(67, 11).Value = 4331.3335  # K67: 4332.3335 -> 4331.3335
$ws.Cells.Item(67, 13).Value = -3473.3335  # M67: -3474.3335 -> -3473.3335
$ws.Cells.Item(69, 8).Value = 12000  # H69: 17500 -> 12000
$ws.Cells.Item(69, 9).Value = 1000  # I69: 0 -> 1000
$ws.Cells.Item(69, 11).Value = 3000  # K69: 0 -> 3000
$ws.Cells.Item(69, 13).Value = -2126  # M69: <MISSING> -> -2126
$ws.Cells.Item(70, 8).Value = 64898.5  # H70: 85031.336 -> 64898.5
$ws.Cells.Item(70, 9).Value = 1899  # I70: 1811.625 -> 1899
$ws.Cells.Item(70, 10).Value = 113898.11  # J70: 251470.75 -> 113898.11
$ws.Cells.Item(70, 11).Value = 5697  # K70: 5434.875 -> 5697
$ws.Cells.Item(70, 12).Value = 341694.33  # L70: 754412.25 -> 341694.33
$ws.Cells.Item(70, 13).Value = -5427  # M70: -5164.875 -> -5427
$ws.Cells.Item(70, 14).Value = -342234.33  # N70: -754952.25 -> -342234.33
$ws.Cells.Item(72, 8).Value = 12000  # H72: 17500 -> 12000
$ws.Cells.Item(72, 9).Value = 1000  # I72: 0 -> 1000
$ws.Cells.Item(72, 11).Value = 9000  # K72: 0 -> 9000
$ws.Cells.Item(72, 13).Value = -4632  # M72: <MISSING> -> -4632
$ws.Cells.Item(73, 8).Value = 64898.5  # H73: 85031.336 -> 64898.5
$ws.Cells.Item(73, 9).Value = 1899  # I73: 1811.625 -> 1899
$ws.Cells.Item(73, 10).Value = 113898.11  # J73: 251470.75 -> 113898.11
$ws.Cells.Item(73, 11).Value = 5697  # K73: 5434.875 -> 5697
$ws.Cells.Item(73, 12).Value = 341694.33  # L73: 754412.25 -> 341694.33
$ws.Cells.Item(73, 13).Value = -4761  # M73: -4498.875 -> -4761
$ws.Cells.Item(73, 14).Value = -343566.33  # N73: -756284.25 -> -343566.33
$ws.Cells.Item(86, 8).Value = 3098125.8  # H86: 3292134.5 -> 3098125.8
$ws.Cells.Item(86, 9).Value = 2218.2222  # I86: 3247.375 -> 2218.2222
$ws.Cells.Item(86, 11).Value = 2218.2222  # K86: 3247.375 -> 2218.2222
$ws.Cells.Item(86, 13).Value = -1095.2222  # M86: -2124.375 -> -1095.2222
$ws.Cells.Item(87, 8).Value = 73550.07000000001  # H87: 66450 -> 73550.07000000001
$ws.Cells.Item(87, 10).Value = 74669.30499999999  # J87: 73900 -> 74669.30499999999
$ws.Cells.Item(87, 12).Value = 74669.30499999999  # L87: 73900 -> 74669.30499999999
$ws.Cells.Item(87, 14).Value = -77165.30499999999  # N87: -76396 -> -77165.30499999999
$ws.Cells.Item(89, 8).Value = 3098125.8  # H89: 3292134.5 -> 3098125.8
$ws.Cells.Item(89, 9).Value = 2218.2222  # I89: 3247.375 -> 2218.2222
$ws.Cells.Item(89, 11).Value = 11091.111  # K89: 16236.875 -> 11091.111
$ws.Cells.Item(89, 13).Value = -5475.111000000001  # M89: -10620.875 -> -5475.111000000001
$ws.Cells.Item(90, 8).Value = 73550.07000000001  # H90: 66450 -> 73550.07000000001
$ws.Cells.Item(90, 10).Value = 74669.30499999999  # J90: 73900 -> 74669.30499999999
$ws.Cells.Item(90, 12).Value = 224007.915  # L90: 221700 -> 224007.915
$ws.Cells.Item(90, 14).Value = -236487.915  # N90: -234180 -> -236487.915
$ws.Cells.Item(92, 8).Value = 151.625  # H92: 162 -> 151.625
$ws.Cells.Item(92, 9).Value = 151.625  # I92: 162 -> 151.625
$ws.Cells.Item(92, 11).Value = 151.625  # K92: 162 -> 151.625
$ws.Cells.Item(92, 13).Value = 1096.375  # M92: 1086 -> 1096.375
$ws.Cells.Item(94, 8).Value = 4792  # H94: 3713.8572 -> 4792
$ws.Cells.Item(94, 9).Value = 6005  # I94: 3511.75 -> 6005
$ws.Cells.Item(94, 11).Value = 6005  # K94: 3511.75 -> 6005
$ws.Cells.Item(94, 13).Value = -5554  # M94: -3060.75 -> -5554
$ws.Cells.Item(96, 8).Value = 509.58334  # H96: 377 -> 509.58334
$ws.Cells.Item(96, 9).Value = 374.72726  # I96: 377 -> 374.72726
$ws.Cells.Item(96, 10).Value = 1993  # J96: 0 -> 1993
$ws.Cells.Item(96, 11).Value = 1124.18178  # K96: 1131 -> 1124.18178
$ws.Cells.Item(96, 12).Value = 5979  # L96: 0 -> 5979
$ws.Cells.Item(96, 13).Value = 248.8182200000001  # M96: 242 -> 248.8182200000001
$ws.Cells.Item(96, 14).Value = -8725  # N96: <MISSING> -> -8725
$ws.Cells.Item(98, 8).Value = 2475.04  # H98: 2505.5806 -> 2475.04
$ws.Cells.Item(98, 9).Value = 1905.7273  # I98: 1859.7826 -> 1905.7273
$ws.Cells.Item(98, 10).Value = 6650  # J98: 4362.25 -> 6650
$ws.Cells.Item(98, 11).Value = 1905.7273  # K98: 1859.7826 -> 1905.7273
$ws.Cells.Item(98, 12).Value = 6650  # L98: 4362.25 -> 6650
$ws.Cells.Item(98, 13).Value = -407.7273  # M98: -361.7826 -> -407.7273
$ws.Cells.Item(98, 14).Value = -9646  # N98: -7358.25 -> -9646
$ws.Cells.Item(100, 8).Value = 1729.6666  # H100: 8088.55 -> 1729.6666
$ws.Cells.Item(100, 9).Value = 1791.7  # I100: 1793.2 -> 1791.7
$ws.Cells.Item(100, 10).Value = 1419.5  # J100: 14383.9 -> 1419.5
$ws.Cells.Item(100, 11).Value = 1791.7  # K100: 1793.2 -> 1791.7
$ws.Cells.Item(100, 12).Value = 1419.5  # L100: 14383.9 -> 1419.5
$ws.Cells.Item(100, 13).Value = -1250.7  # M100: -1252.2 -> -1250.7
$ws.Cells.Item(100, 14).Value = -2501.5  # N100: -15465.9 -> -2501.5
$ws.Cells.Item(107, 8).Value = 56483.11  # H107: 53612.95 -> 56483.11
$ws.Cells.Item(107, 10).Value = 0  # J107: 1950 -> 0
$ws.Cells.Item(107, 12).Value = 0  # L107: 1950 -> 0
$ws.Cells.Item(107, 14).ClearContents()  # N107: -5790 -> (cleared)
$ws.Cells.Item(108, 8).Value = 48969.332  # H108: 48962 -> 48969.332
$ws.Cells.Item(108, 10).Value = 48969.332  # J108: 48962 -> 48969.332
$ws.Cells.Item(108, 12).Value = 48969.332  # L108: 48962 -> 48969.332
$ws.Cells.Item(108, 14).Value = -56649.332  # N108: -56642 -> -56649.332
$ws.Cells.Item(111, 8).Value = 104744.4  # H111: 95361.09 -> 104744.4
$ws.Cells.Item(111, 9).Value = 171071.33  # I111: 146850.86 -> 171071.33
$ws.Cells.Item(111, 11).Value = 513213.99  # K111: 440552.58 -> 513213.99
$ws.Cells.Item(111, 13).Value = -510146.99  # M111: -437485.58 -> -510146.99
$ws.Cells.Item(112, 8).Value = 5079.143  # H112: 4865.1333 -> 5079.143
$ws.Cells.Item(112, 10).Value = 5079.143  # J112: 4865.1333 -> 5079.143
$ws.Cells.Item(112, 12).Value = 15237.429  # L112: 14595.3999 -> 15237.429
$ws.Cells.Item(112, 14).Value = -17453.429  # N112: -16811.3999 -> -17453.429
$ws.Cells.Item(116, 8).Value = 6750  # H116: 7311.875 -> 6750
$ws.Cells.Item(116, 10).Value = 0  # J116: 8997.5 -> 0
$ws.Cells.Item(116, 12).Value = 0  # L116: 8997.5 -> 0
$ws.Cells.Item(116, 14).ClearContents()  # N116: -15881.5 -> (cleared)
$ws.Cells.Item(118, 8).Value = 635.8  # H118: 779 -> 635.8
$ws.Cells.Item(118, 9).Value = 528.6667  # I118: 571.125 -> 528.6667
$ws.Cells.Item(118, 10).Value = 1600  # J118: 1333.3334 -> 1600
$ws.Cells.Item(118, 11).Value = 1586.0001  # K118: 1713.375 -> 1586.0001
$ws.Cells.Item(118, 12).Value = 4800  # L118: 4000.0002 -> 4800
$ws.Cells.Item(118, 13).Value = 70.99990000000003  # M118: -56.375 -> 70.99990000000003
$ws.Cells.Item(118, 14).Value = -8114  # N118: -7314.0002 -> -8114
$ws.Cells.Item(122, 8).Value = 2475.04  # H122: 2505.5806 -> 2475.04
$ws.Cells.Item(122, 9).Value = 1905.7273  # I122: 1859.7826 -> 1905.7273
$ws.Cells.Item(122, 10).Value = 6650  # J122: 4362.25 -> 6650
$ws.Cells.Item(122, 11).Value = 5717.1819  # K122: 5579.3478 -> 5717.1819
$ws.Cells.Item(122, 12).Value = 19950  # L122: 13086.75 -> 19950
$ws.Cells.Item(122, 13).Value = -3267.1819  # M122: -3129.3478 -> -3267.1819
$ws.Cells.Item(122, 14).Value = -24850  # N122: -17986.75 -> -24850
$ws.Cells.Item(125, 8).Value = 2754.6667  # H125: 2648.3157 -> 2754.6667
$ws.Cells.Item(125, 10).Value = 5128  # J125: 4578.75 -> 5128
$ws.Cells.Item(125, 12).Value = 46152  # L125: 41208.75 -> 46152
$ws.Cells.Item(125, 14).Value = -51072  # N125: -46128.75 -> -51072
$ws.Cells.Item(127, 8).Value = 14238  # H127: 11804.637 -> 14238
$ws.Cells.Item(127, 9).Value = 18837.666  # I127: 16207.571 -> 18837.666
$ws.Cells.Item(127, 10).Value = 5038.6665  # J127: 4099.5 -> 5038.6665
$ws.Cells.Item(127, 11).Value = 56512.99800000001  # K127: 48622.713 -> 56512.99800000001
$ws.Cells.Item(127, 12).Value = 15115.9995  # L127: 12298.5 -> 15115.9995
$ws.Cells.Item(127, 13).Value = -51552.99800000001  # M127: -43662.713 -> -51552.99800000001
$ws.Cells.Item(127, 14).Value = -25035.9995  # N127: -22218.5 -> -25035.9995
$ws.Cells.Item(129, 8).Value = 7655.436  # H129: 7809.6055 -> 7655.436
$ws.Cells.Item(129, 9).Value = 12984  # I129: 12997.333 -> 12984
$ws.Cells.Item(129, 10).Value = 6056.8667  # J129: 6199.6206 -> 6056.8667
$ws.Cells.Item(129, 11).Value = 38952  # K129: 38991.999 -> 38952
$ws.Cells.Item(129, 12).Value = 18170.6001  # L129: 18598.8618 -> 18170.6001
$ws.Cells.Item(129, 13).Value = -33952  # M129: -33991.999 -> -33952
$ws.Cells.Item(129, 14).Value = -28170.6001  # N129: -28598.8618 -> -28170.6001
$ws.Cells.Item(131, 8).Value = 1448.3  # H131: 2960.84 -> 1448.3
$ws.Cells.Item(131, 9).Value = 1448.3  # I131: 1527.421 -> 1448.3
$ws.Cells.Item(131, 10).Value = 0  # J131: 7500 -> 0
$ws.Cells.Item(131, 11).Value = 4344.9  # K131: 4582.263 -> 4344.9
$ws.Cells.Item(131, 12).Value = 0  # L131: 22500 -> 0
$ws.Cells.Item(131, 13).Value = 695.1000000000004  # M131: 457.7370000000001 -> 695.1000000000004
$ws.Cells.Item(131, 14).ClearContents()  # N131: -32580 -> (cleared)
$ws.Cells.Item(132, 8).Value = 2093.1316  # H132: 2162.0273 -> 2093.1316
$ws.Cells.Item(132, 9).Value = 1642.8616  # I132: 1702.1936 -> 1642.8616
$ws.Cells.Item(132, 11).Value = 4928.5848  # K132: 5106.5808 -> 4928.5848
$ws.Cells.Item(132, 13).Value = -2398.5848  # M132: -2576.5808 -> -2398.5848
$ws.Cells.Item(135, 8).Value = 541764.7  # H135: 572698.3 -> 541764.7
$ws.Cells.Item(135, 9).Value = 770256.4  # I135: 834408.8 -> 770256.4
$ws.Cells.Item(135, 11).Value = 6932307.600000001  # K135: 7509679.2 -> 6932307.600000001
$ws.Cells.Item(135, 13).Value = -6929772.600000001  # M135: -7507144.2 -> -6929772.600000001
$ws.Cells.Item(137, 8).Value = 649174.3  # H137: 693858.75 -> 649174.3
$ws.Cells.Item(137, 9).Value = 529073.6  # I137: 558394.4 -> 529073.6
$ws.Cells.Item(137, 10).Value = 839333.75  # J137: 915527.75 -> 839333.75
$ws.Cells.Item(137, 11).Value = 1587220.8  # K137: 1675183.2 -> 1587220.8
$ws.Cells.Item(137, 12).Value = 2518001.25  # L137: 2746583.25 -> 2518001.25
$ws.Cells.Item(137, 13).Value = -1584670.8  # M137: -1672633.2 -> -1584670.8
$ws.Cells.Item(137, 14).Value = -2523101.25  # N137: -2751683.25 -> -2523101.25
$ws.Cells.Item(138, 8).Value = 3960.746  # H138: 3986.238 -> 3960.746
$ws.Cells.Item(138, 9).Value = 2041.1482  # I138: 2068.1482 -> 2041.1482
$ws.Cells.Item(138, 10).Value = 5400.4443  # J138: 5424.8057 -> 5400.4443
$ws.Cells.Item(138, 11).Value = 6123.444600000001  # K138: 6204.444600000001 -> 6123.444600000001
$ws.Cells.Item(138, 12).Value = 16201.3329  # L138: 16274.4171 -> 16201.3329
$ws.Cells.Item(138, 13).Value = -983.4446000000007  # M138: -1064.444600000001 -> -983.4446000000007
$ws.Cells.Item(138, 14).Value = -26481.3329  # N138: -26554.4171 -> -26481.3329
$ws.Cells.Item(141, 8).Value = 2127.7334  # H141: 2098.1147 -> 2127.7334
$ws.Cells.Item(141, 9).Value = 958.72  # I141: 970.67346 -> 958.72
$ws.Cells.Item(141, 10).Value = 7972.8  # J141: 6701.8335 -> 7972.8
$ws.Cells.Item(141, 11).Value = 2876.16  # K141: 2912.02038 -> 2876.16
$ws.Cells.Item(141, 12).Value = 23918.4  # L141: 20105.5005 -> 23918.4
$ws.Cells.Item(141, 13).Value = 2303.84  # M141: 2267.97962 -> 2303.84
$ws.Cells.Item(141, 14).Value = -34278.4  # N141: -30465.5005 -> -34278.4

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 59300.316  # H2: 62811.39 -> 59300.316
$ws.Cells.Item(2, 9).Value = 62317.11  # I2: 70038 -> 62317.11
$ws.Cells.Item(2, 10).Value = 4998  # J2: 4998.5 -> 4998
$ws.Cells.Item(2, 11).Value = 62317.11  # K2: 70038 -> 62317.11
$ws.Cells.Item(2, 12).Value = 4998  # L2: 4998.5 -> 4998
$ws.Cells.Item(2, 13).Value = -62204.11  # M2: -69925 -> -62204.11
$ws.Cells.Item(2, 14).Value = -5224  # N2: -5224.5 -> -5224
$ws.Cells.Item(32, 8).Value = 7147.4194  # H32: 8073.3706 -> 7147.4194
$ws.Cells.Item(32, 9).Value = 6818.077  # I32: 7586.1304 -> 6818.077
$ws.Cells.Item(32, 10).Value = 8860  # J32: 10875 -> 8860
$ws.Cells.Item(32, 11).Value = 6818.077  # K32: 7586.1304 -> 6818.077
$ws.Cells.Item(32, 12).Value = 8860  # L32: 10875 -> 8860
$ws.Cells.Item(32, 13).Value = -6531.077  # M32: -7299.1304 -> -6531.077
$ws.Cells.Item(32, 14).Value = -9434  # N32: -11449 -> -9434
$ws.Cells.Item(45, 8).Value = 3650.5386  # H45: 3682 -> 3650.5386
$ws.Cells.Item(45, 9).Value = 2550.7778  # I45: 2645.875 -> 2550.7778
$ws.Cells.Item(45, 10).Value = 6125  # J45: 5339.8 -> 6125
$ws.Cells.Item(45, 11).Value = 2550.7778  # K45: 2645.875 -> 2550.7778
$ws.Cells.Item(45, 12).Value = 6125  # L45: 5339.8 -> 6125
$ws.Cells.Item(45, 13).Value = -2173.7778  # M45: -2268.875 -> -2173.7778
$ws.Cells.Item(45, 14).Value = -6879  # N45: -6093.8 -> -6879
$ws.Cells.Item(61, 8).Value = 1828.7368  # H61: 1798.7693 -> 1828.7368
$ws.Cells.Item(61, 9).Value = 931.5  # I61: 923.2727 -> 931.5
$ws.Cells.Item(61, 11).Value = 931.5  # K61: 923.2727 -> 931.5
$ws.Cells.Item(61, 13).Value = -719.5  # M61: -711.2727 -> -719.5
$ws.Cells.Item(88, 8).Value = 2511.0527  # H88: 2831.875 -> 2511.0527
$ws.Cells.Item(88, 9).Value = 2668.4546  # I88: 2925.3 -> 2668.4546
$ws.Cells.Item(88, 10).Value = 2294.625  # J88: 2676.1667 -> 2294.625
$ws.Cells.Item(88, 11).Value = 2668.4546  # K88: 2925.3 -> 2668.4546
$ws.Cells.Item(88, 12).Value = 2294.625  # L88: 2676.1667 -> 2294.625
$ws.Cells.Item(88, 13).Value = -2262.4546  # M88: -2519.3 -> -2262.4546
$ws.Cells.Item(88, 14).Value = -3106.625  # N88: -3488.1667 -> -3106.625
$ws.Cells.Item(91, 8).Value = 2511.0527  # H91: 2831.875 -> 2511.0527
$ws.Cells.Item(91, 9).Value = 2668.4546  # I91: 2925.3 -> 2668.4546
$ws.Cells.Item(91, 10).Value = 2294.625  # J91: 2676.1667 -> 2294.625
$ws.Cells.Item(91, 11).Value = 2668.4546  # K91: 2925.3 -> 2668.4546
$ws.Cells.Item(91, 12).Value = 2294.625  # L91: 2676.1667 -> 2294.625
$ws.Cells.Item(91, 13).Value = -1264.4546  # M91: -1521.3 -> -1264.4546
$ws.Cells.Item(91, 14).Value = -5102.625  # N91: -5484.1667 -> -5102.625
$ws.Cells.Item(97, 8).Value = 3975.8147  # H97: 4791.115 -> 3975.8147
$ws.Cells.Item(97, 9).Value = 4405.8184  # I97: 4807.55 -> 4405.8184
$ws.Cells.Item(97, 10).Value = 2083.8  # J97: 4736.3335 -> 2083.8
$ws.Cells.Item(97, 11).Value = 4405.8184  # K97: 4807.55 -> 4405.8184
$ws.Cells.Item(97, 12).Value = 2083.8  # L97: 4736.3335 -> 2083.8
$ws.Cells.Item(97, 13).Value = -3909.8184  # M97: -4311.55 -> -3909.8184
$ws.Cells.Item(97, 14).Value = -3075.8  # N97: -5728.3335 -> -3075.8
$ws.Cells.Item(102, 8).Value = 3324.75  # H102: 5600 -> 3324.75
$ws.Cells.Item(102, 9).Value = 3324.75  # I102: 5600 -> 3324.75
$ws.Cells.Item(102, 11).Value = 3324.75  # K102: 5600 -> 3324.75
$ws.Cells.Item(102, 13).Value = -1702.75  # M102: -3978 -> -1702.75
$ws.Cells.Item(110, 8).Value = 278960.38  # H110: 264160.78 -> 278960.38
$ws.Cells.Item(110, 9).Value = 313455.75  # I110: 278669.84 -> 313455.75
$ws.Cells.Item(110, 10).Value = 2997.5  # J110: 2998 -> 2997.5
$ws.Cells.Item(110, 11).Value = 313455.75  # K110: 278669.84 -> 313455.75
$ws.Cells.Item(110, 12).Value = 2997.5  # L110: 2998 -> 2997.5
$ws.Cells.Item(110, 13).Value = -311410.75  # M110: -276624.84 -> -311410.75
$ws.Cells.Item(110, 14).Value = -7087.5  # N110: -7088 -> -7087.5
$ws.Cells.Item(116, 8).Value = 59300.316  # H116: 62811.39 -> 59300.316
$ws.Cells.Item(116, 9).Value = 62317.11  # I116: 70038 -> 62317.11
$ws.Cells.Item(116, 10).Value = 4998  # J116: 4998.5 -> 4998
$ws.Cells.Item(116, 11).Value = 62317.11  # K116: 70038 -> 62317.11
$ws.Cells.Item(116, 12).Value = 4998  # L116: 4998.5 -> 4998
$ws.Cells.Item(116, 13).Value = -60023.11  # M116: -67744 -> -60023.11
$ws.Cells.Item(116, 14).Value = -9586  # N116: -9586.5 -> -9586
$ws.Cells.Item(132, 8).Value = 3826.8462  # H132: 3982.2163 -> 3826.8462
$ws.Cells.Item(132, 9).Value = 2532.1936  # I132: 2641.138 -> 2532.1936
$ws.Cells.Item(132, 11).Value = 7596.5808  # K132: 7923.414 -> 7596.5808
$ws.Cells.Item(132, 13).Value = -5066.5808  # M132: -5393.414 -> -5066.5808
$ws.Cells.Item(136, 8).Value = 1828.7368  # H136: 1798.7693 -> 1828.7368
$ws.Cells.Item(136, 9).Value = 931.5  # I136: 923.2727 -> 931.5
$ws.Cells.Item(136, 11).Value = 2794.5  # K136: 2769.8181 -> 2794.5
$ws.Cells.Item(136, 13).Value = -244.5  # M136: -219.8181 -> -244.5

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 59300.316  # H3: 62811.39 -> 59300.316
$ws.Cells.Item(3, 9).Value = 62317.11  # I3: 70038 -> 62317.11
$ws.Cells.Item(3, 10).Value = 4998  # J3: 4998.5 -> 4998
$ws.Cells.Item(3, 11).Value = 62317.11  # K3: 70038 -> 62317.11
$ws.Cells.Item(3, 12).Value = 4998  # L3: 4998.5 -> 4998
$ws.Cells.Item(3, 13).Value = -62203.11  # M3: -69924 -> -62203.11
$ws.Cells.Item(3, 14).Value = -5226  # N3: -5226.5 -> -5226
$ws.Cells.Item(86, 8).Value = 1135225.4  # H86: 1309722.4 -> 1135225.4
$ws.Cells.Item(86, 9).Value = 1216170  # I86: 1418699.2 -> 1216170
$ws.Cells.Item(86, 11).Value = 1216170  # K86: 1418699.2 -> 1216170
$ws.Cells.Item(86, 13).Value = -1215047  # M86: -1417576.2 -> -1215047
$ws.Cells.Item(89, 8).Value = 1135225.4  # H89: 1309722.4 -> 1135225.4
$ws.Cells.Item(89, 9).Value = 1216170  # I89: 1418699.2 -> 1216170
$ws.Cells.Item(89, 11).Value = 6080850  # K89: 7093496 -> 6080850
$ws.Cells.Item(89, 13).Value = -6075234  # M89: -7087880 -> -6075234
$ws.Cells.Item(94, 8).Value = 4804.643  # H94: 5113.846 -> 4804.643
$ws.Cells.Item(94, 9).Value = 2251.4546  # I94: 2398.1 -> 2251.4546
$ws.Cells.Item(94, 11).Value = 2251.4546  # K94: 2398.1 -> 2251.4546
$ws.Cells.Item(94, 13).Value = -1800.4546  # M94: -1947.1 -> -1800.4546
$ws.Cells.Item(105, 8).Value = 64843.812  # H105: 68970.13 -> 64843.812
$ws.Cells.Item(105, 9).Value = 79372.30499999999  # I105: 73752.92999999999 -> 79372.30499999999
$ws.Cells.Item(105, 10).Value = 1887  # J105: 2011 -> 1887
$ws.Cells.Item(105, 11).Value = 79372.30499999999  # K105: 73752.92999999999 -> 79372.30499999999
$ws.Cells.Item(105, 12).Value = 1887  # L105: 2011 -> 1887
$ws.Cells.Item(105, 13).Value = -77625.30499999999  # M105: -72005.92999999999 -> -77625.30499999999
$ws.Cells.Item(105, 14).Value = -5381  # N105: -5505 -> -5381
$ws.Cells.Item(107, 8).Value = 186887.89  # H107: 186888.38 -> 186887.89
$ws.Cells.Item(107, 9).Value = 1576.738  # I107: 1577.3572 -> 1576.738
$ws.Cells.Item(107, 11).Value = 1576.738  # K107: 1577.3572 -> 1576.738
$ws.Cells.Item(107, 13).Value = 343.2619999999999  # M107: 342.6428000000001 -> 343.2619999999999

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 330.6389  # H7: 333.75 -> 330.6389
$ws.Cells.Item(7, 9).Value = 247.14815  # I7: 251.2963 -> 247.14815
$ws.Cells.Item(7, 11).Value = 247.14815  # K7: 251.2963 -> 247.14815
$ws.Cells.Item(7, 13).Value = -134.14815  # M7: -138.2963 -> -134.14815
$ws.Cells.Item(16, 8).Value = 3678.9092  # H16: 3747.5 -> 3678.9092
$ws.Cells.Item(16, 10).Value = 2996.4  # J16: 2997.25 -> 2996.4
$ws.Cells.Item(16, 12).Value = 2996.4  # L16: 2997.25 -> 2996.4
$ws.Cells.Item(16, 14).Value = -3570.4  # N16: -3571.25 -> -3570.4
$ws.Cells.Item(22, 8).Value = 794.61536  # H22: 827.5 -> 794.61536
$ws.Cells.Item(22, 9).Value = 794.61536  # I22: 827.5 -> 794.61536
$ws.Cells.Item(22, 11).Value = 794.61536  # K22: 827.5 -> 794.61536
$ws.Cells.Item(22, 13).Value = -444.61536  # M22: -477.5 -> -444.61536
$ws.Cells.Item(31, 8).Value = 527589  # H31: 527784.75 -> 527589
$ws.Cells.Item(31, 9).Value = 304779.8  # I31: 304889.75 -> 304779.8
$ws.Cells.Item(31, 10).Value = 1577975.2  # J31: 1578575.2 -> 1577975.2
$ws.Cells.Item(31, 11).Value = 304779.8  # K31: 304889.75 -> 304779.8
$ws.Cells.Item(31, 12).Value = 1577975.2  # L31: 1578575.2 -> 1577975.2
$ws.Cells.Item(31, 13).Value = -304484.8  # M31: -304594.75 -> -304484.8
$ws.Cells.Item(31, 14).Value = -1578565.2  # N31: -1579165.2 -> -1578565.2
$ws.Cells.Item(34, 8).Value = 527589  # H34: 527784.75 -> 527589
$ws.Cells.Item(34, 9).Value = 304779.8  # I34: 304889.75 -> 304779.8
$ws.Cells.Item(34, 10).Value = 1577975.2  # J34: 1578575.2 -> 1577975.2
$ws.Cells.Item(34, 11).Value = 304779.8  # K34: 304889.75 -> 304779.8
$ws.Cells.Item(34, 12).Value = 1577975.2  # L34: 1578575.2 -> 1577975.2
$ws.Cells.Item(34, 13).Value = -304577.8  # M34: -304687.75 -> -304577.8
$ws.Cells.Item(34, 14).Value = -1578379.2  # N34: -1578979.2 -> -1578379.2
$ws.Cells.Item(58, 8).Value = 203926.45  # H58: 212388.08 -> 203926.45
$ws.Cells.Item(58, 9).Value = 288145.22  # I58: 305557.22 -> 288145.22
$ws.Cells.Item(58, 11).Value = 288145.22  # K58: 305557.22 -> 288145.22
$ws.Cells.Item(58, 13).Value = -287942.22  # M58: -305354.22 -> -287942.22
$ws.Cells.Item(62, 8).Value = 2792.625  # H62: 3306.9167 -> 2792.625
$ws.Cells.Item(62, 9).Value = 2181.5454  # I62: 2666.3333 -> 2181.5454
$ws.Cells.Item(62, 10).Value = 4137  # J62: 3947.5 -> 4137
$ws.Cells.Item(62, 11).Value = 2181.5454  # K62: 2666.3333 -> 2181.5454
$ws.Cells.Item(62, 12).Value = 4137  # L62: 3947.5 -> 4137
$ws.Cells.Item(62, 13).Value = -1557.5454  # M62: -2042.3333 -> -1557.5454
$ws.Cells.Item(62, 14).Value = -5385  # N62: -5195.5 -> -5385
$ws.Cells.Item(65, 8).Value = 2792.625  # H65: 3306.9167 -> 2792.625
$ws.Cells.Item(65, 9).Value = 2181.5454  # I65: 2666.3333 -> 2181.5454
$ws.Cells.Item(65, 10).Value = 4137  # J65: 3947.5 -> 4137
$ws.Cells.Item(65, 11).Value = 10907.727  # K65: 13331.6665 -> 10907.727
$ws.Cells.Item(65, 12).Value = 20685  # L65: 19737.5 -> 20685
$ws.Cells.Item(65, 13).Value = -7787.726999999999  # M65: -10211.6665 -> -7787.726999999999
$ws.Cells.Item(65, 14).Value = -26925  # N65: -25977.5 -> -26925
$ws.Cells.Item(81, 8).Value = 0  # H81: 140000 -> 0
$ws.Cells.Item(81, 10).Value = 0  # J81: 140000 -> 0
$ws.Cells.Item(81, 12).Value = 0  # L81: 140000 -> 0
$ws.Cells.Item(81, 14).ClearContents()  # N81: -141996 -> (cleared)
$ws.Cells.Item(84, 8).Value = 0  # H84: 140000 -> 0
$ws.Cells.Item(84, 10).Value = 0  # J84: 140000 -> 0
$ws.Cells.Item(84, 12).Value = 0  # L84: 420000 -> 0
$ws.Cells.Item(84, 14).ClearContents()  # N84: -429984 -> (cleared)
$ws.Cells.Item(86, 8).Value = 11890  # H86: 12223.556 -> 11890
$ws.Cells.Item(86, 10).Value = 12002  # J86: 13002.667 -> 12002
$ws.Cells.Item(86, 12).Value = 12002  # L86: 13002.667 -> 12002
$ws.Cells.Item(86, 14).Value = -14248  # N86: -15248.667 -> -14248
$ws.Cells.Item(89, 8).Value = 11890  # H89: 12223.556 -> 11890
$ws.Cells.Item(89, 10).Value = 12002  # J89: 13002.667 -> 12002
$ws.Cells.Item(89, 12).Value = 60010  # L89: 65013.335 -> 60010
$ws.Cells.Item(89, 14).Value = -71242  # N89: -76245.33499999999 -> -71242
$ws.Cells.Item(99, 8).Value = 5249.3  # H99: 5371.724 -> 5249.3
$ws.Cells.Item(99, 9).Value = 4201.381  # I99: 4326.5 -> 4201.381
$ws.Cells.Item(99, 11).Value = 4201.381  # K99: 4326.5 -> 4201.381
$ws.Cells.Item(99, 13).Value = -2703.381  # M99: -2828.5 -> -2703.381
$ws.Cells.Item(105, 8).Value = 1852.6666  # H105: 1817.2632 -> 1852.6666
$ws.Cells.Item(105, 9).Value = 1770.2222  # I105: 1711.2 -> 1770.2222
$ws.Cells.Item(105, 11).Value = 1770.2222  # K105: 1711.2 -> 1770.2222
$ws.Cells.Item(105, 13).Value = -23.22219999999993  # M105: 35.79999999999995 -> -23.22219999999993
$ws.Cells.Item(113, 8).Value = 3678.9092  # H113: 3747.5 -> 3678.9092
$ws.Cells.Item(113, 10).Value = 2996.4  # J113: 2997.25 -> 2996.4
$ws.Cells.Item(113, 12).Value = 2996.4  # L113: 2997.25 -> 2996.4
$ws.Cells.Item(113, 14).Value = -7336.4  # N113: -7337.25 -> -7336.4
$ws.Cells.Item(126, 8).Value = 5249.3  # H126: 5371.724 -> 5249.3
$ws.Cells.Item(126, 9).Value = 4201.381  # I126: 4326.5 -> 4201.381
$ws.Cells.Item(126, 11).Value = 12604.143  # K126: 12979.5 -> 12604.143
$ws.Cells.Item(126, 13).Value = -10134.143  # M126: -10509.5 -> -10134.143
$ws.Cells.Item(132, 8).Value = 2807.2068  # H132: 2292.6184 -> 2807.2068
$ws.Cells.Item(132, 9).Value = 2108.9148  # I132: 1720.9375 -> 2108.9148
$ws.Cells.Item(132, 10).Value = 5790.8184  # J132: 5341.5835 -> 5790.8184
$ws.Cells.Item(132, 11).Value = 6326.7444  # K132: 5162.8125 -> 6326.7444
$ws.Cells.Item(132, 12).Value = 17372.4552  # L132: 16024.7505 -> 17372.4552
$ws.Cells.Item(132, 13).Value = -3796.7444  # M132: -2632.8125 -> -3796.7444
$ws.Cells.Item(132, 14).Value = -22432.4552  # N132: -21084.7505 -> -22432.4552
$ws.Cells.Item(134, 8).Value = 330906.25  # H134: 330907.38 -> 330906.25
$ws.Cells.Item(134, 9).Value = 202378.77  # I134: 202380.1 -> 202378.77
$ws.Cells.Item(134, 11).Value = 607136.3099999999  # K134: 607140.3 -> 607136.3099999999
$ws.Cells.Item(134, 13).Value = -604601.3099999999  # M134: -604605.3 -> -604601.3099999999
$ws.Cells.Item(136, 8).Value = 203926.45  # H136: 212388.08 -> 203926.45
$ws.Cells.Item(136, 9).Value = 288145.22  # I136: 305557.22 -> 288145.22
$ws.Cells.Item(136, 11).Value = 864435.6599999999  # K136: 916671.6599999999 -> 864435.6599999999
$ws.Cells.Item(136, 13).Value = -861885.6599999999  # M136: -914121.6599999999 -> -861885.6599999999

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 608.94446  # H12: 174.21428 -> 608.94446
$ws.Cells.Item(12, 9).Value = 102.25  # I12: 200 -> 102.25
$ws.Cells.Item(12, 10).Value = 753.7143  # J12: 172.23077 -> 753.7143
$ws.Cells.Item(12, 11).Value = 306.75  # K12: 600 -> 306.75
$ws.Cells.Item(12, 12).Value = 2261.1429  # L12: 516.69231 -> 2261.1429
$ws.Cells.Item(12, 13).Value = -133.75  # M12: -427 -> -133.75
$ws.Cells.Item(12, 14).Value = -2607.1429  # N12: -862.69231 -> -2607.1429
$ws.Cells.Item(17, 8).Value = 1266  # H17: 1299.3334 -> 1266
$ws.Cells.Item(17, 9).Value = 1266  # I17: 1299 -> 1266
$ws.Cells.Item(17, 10).Value = 0  # J17: 1300 -> 0
$ws.Cells.Item(17, 11).Value = 3798  # K17: 3897 -> 3798
$ws.Cells.Item(17, 12).Value = 0  # L17: 3900 -> 0
$ws.Cells.Item(17, 13).Value = -3629  # M17: -3728 -> -3629
$ws.Cells.Item(17, 14).ClearContents()  # N17: -4238 -> (cleared)
$ws.Cells.Item(23, 8).Value = 427.83334  # H23: 439.7647 -> 427.83334
$ws.Cells.Item(23, 9).Value = 350.5  # I23: 375.6 -> 350.5
$ws.Cells.Item(23, 11).Value = 1051.5  # K23: 1126.8 -> 1051.5
$ws.Cells.Item(23, 13).Value = -816.5  # M23: -891.8000000000002 -> -816.5
$ws.Cells.Item(57, 8).Value = 4375  # H57: 4500 -> 4375
$ws.Cells.Item(57, 10).Value = 0  # J57: 5000 -> 0
$ws.Cells.Item(57, 12).Value = 0  # L57: 15000 -> 0
$ws.Cells.Item(57, 14).ClearContents()  # N57: -16118 -> (cleared)
$ws.Cells.Item(129, 8).Value = 67217.53  # H129: 67224.266 -> 67217.53
$ws.Cells.Item(129, 9).Value = 329.84616  # I129: 353.25 -> 329.84616
$ws.Cells.Item(129, 10).Value = 501987.5  # J129: 334708.34 -> 501987.5
$ws.Cells.Item(129, 11).Value = 989.5384799999999  # K129: 1059.75 -> 989.5384799999999
$ws.Cells.Item(129, 12).Value = 1505962.5  # L129: 1004125.02 -> 1505962.5
$ws.Cells.Item(129, 13).Value = 4010.46152  # M129: 3940.25 -> 4010.46152
$ws.Cells.Item(129, 14).Value = -1515962.5  # N129: -1014125.02 -> -1515962.5
$ws.Cells.Item(131, 8).Value = 2910.5933  # H131: 2970.9473 -> 2910.5933
$ws.Cells.Item(131, 9).Value = 1218.2222  # I131: 1308 -> 1218.2222
$ws.Cells.Item(131, 10).Value = 3215.22  # J131: 3242.449 -> 3215.22
$ws.Cells.Item(131, 11).Value = 3654.6666  # K131: 3924 -> 3654.6666
$ws.Cells.Item(131, 12).Value = 9645.66  # L131: 9727.347 -> 9645.66
$ws.Cells.Item(131, 13).Value = 1385.3334  # M131: 1116 -> 1385.3334
$ws.Cells.Item(131, 14).Value = -19725.66  # N131: -19807.347 -> -19725.66
$ws.Cells.Item(132, 8).Value = 849914  # H132: 736981.8 -> 849914
$ws.Cells.Item(132, 9).Value = 250871.75  # I132: 200897.4 -> 250871.75
$ws.Cells.Item(132, 10).Value = 1116155  # J132: 1005024 -> 1116155
$ws.Cells.Item(132, 11).Value = 2257845.75  # K132: 1808076.6 -> 2257845.75
$ws.Cells.Item(132, 12).Value = 10045395  # L132: 9045216 -> 10045395
$ws.Cells.Item(132, 13).Value = -2255315.75  # M132: -1805546.6 -> -2255315.75
$ws.Cells.Item(132, 14).Value = -10050455  # N132: -9050276 -> -10050455
$ws.Cells.Item(136, 8).Value = 7848.5  # H136: 7741.75 -> 7848.5
$ws.Cells.Item(136, 9).Value = 6072.4287  # I136: 6324.6665 -> 6072.4287
$ws.Cells.Item(136, 10).Value = 11992.667  # J136: 11993 -> 11992.667
$ws.Cells.Item(136, 11).Value = 18217.2861  # K136: 18973.9995 -> 18217.2861
$ws.Cells.Item(136, 12).Value = 35978.001  # L136: 35979 -> 35978.001
$ws.Cells.Item(136, 13).Value = -13117.2861  # M136: -13873.9995 -> -13117.2861
$ws.Cells.Item(136, 14).Value = -46178.001  # N136: -46179 -> -46178.001
$ws.Cells.Item(137, 8).Value = 3971.5386  # H137: 4063.8462 -> 3971.5386
$ws.Cells.Item(137, 9).Value = 1856.7142  # I137: 2124.625 -> 1856.7142
$ws.Cells.Item(137, 10).Value = 6438.8335  # J137: 7166.6 -> 6438.8335
$ws.Cells.Item(137, 11).Value = 5570.142599999999  # K137: 6373.875 -> 5570.142599999999
$ws.Cells.Item(137, 12).Value = 19316.5005  # L137: 21499.8 -> 19316.5005
$ws.Cells.Item(137, 13).Value = -470.1425999999992  # M137: -1273.875 -> -470.1425999999992
$ws.Cells.Item(137, 14).Value = -29516.5005  # N137: -31699.8 -> -29516.5005

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 108.9375  # H2: 109 -> 108.9375
$ws.Cells.Item(2, 9).Value = 121.454544  # I2: 121.545456 -> 121.454544
$ws.Cells.Item(2, 11).Value = 121.454544  # K2: 121.545456 -> 121.454544
$ws.Cells.Item(2, 13).Value = -8.454543999999999  # M2: -8.545456000000001 -> -8.454543999999999
$ws.Cells.Item(80, 8).Value = 1254907.5  # H80: 1672750.9 -> 1254907.5
$ws.Cells.Item(80, 9).Value = 1254502.4  # I80: 2006601.8 -> 1254502.4
$ws.Cells.Item(80, 10).Value = 1255312.8  # J80: 1434286 -> 1255312.8
$ws.Cells.Item(80, 11).Value = 1254502.4  # K80: 2006601.8 -> 1254502.4
$ws.Cells.Item(80, 12).Value = 1255312.8  # L80: 1434286 -> 1255312.8
$ws.Cells.Item(80, 13).Value = -1253504.4  # M80: -2005603.8 -> -1253504.4
$ws.Cells.Item(80, 14).Value = -1257308.8  # N80: -1436282 -> -1257308.8
$ws.Cells.Item(83, 8).Value = 1254907.5  # H83: 1672750.9 -> 1254907.5
$ws.Cells.Item(83, 9).Value = 1254502.4  # I83: 2006601.8 -> 1254502.4
$ws.Cells.Item(83, 10).Value = 1255312.8  # J83: 1434286 -> 1255312.8
$ws.Cells.Item(83, 11).Value = 6272512  # K83: 10033009 -> 6272512
$ws.Cells.Item(83, 12).Value = 6276564  # L83: 7171430 -> 6276564
$ws.Cells.Item(83, 13).Value = -6267520  # M83: -10028017 -> -6267520
$ws.Cells.Item(83, 14).Value = -6286548  # N83: -7181414 -> -6286548
$ws.Cells.Item(102, 8).Value = 1693.1837  # H102: 1721.5 -> 1693.1837
$ws.Cells.Item(102, 9).Value = 934.80554  # I102: 935.25 -> 934.80554
$ws.Cells.Item(102, 10).Value = 3793.3076  # J102: 4080.25 -> 3793.3076
$ws.Cells.Item(102, 11).Value = 934.80554  # K102: 935.25 -> 934.80554
$ws.Cells.Item(102, 12).Value = 3793.3076  # L102: 4080.25 -> 3793.3076
$ws.Cells.Item(102, 13).Value = 687.19446  # M102: 686.75 -> 687.19446
$ws.Cells.Item(102, 14).Value = -7037.3076  # N102: -7324.25 -> -7037.3076
$ws.Cells.Item(122, 8).Value = 4277.5454  # H122: 4500.4443 -> 4277.5454
$ws.Cells.Item(122, 9).Value = 3506.625  # I122: 3584 -> 3506.625
$ws.Cells.Item(122, 11).Value = 10519.875  # K122: 10752 -> 10519.875
$ws.Cells.Item(122, 13).Value = -8069.875  # M122: -8302 -> -8069.875
$ws.Cells.Item(126, 8).Value = 3812.8635  # H126: 4726.25 -> 3812.8635
$ws.Cells.Item(126, 10).Value = 4261.2  # J126: 7350 -> 4261.2
$ws.Cells.Item(126, 12).Value = 12783.6  # L126: 22050 -> 12783.6
$ws.Cells.Item(126, 14).Value = -17723.6  # N126: -26990 -> -17723.6
$ws.Cells.Item(132, 8).Value = 195620.83  # H132: 206470.19 -> 195620.83
$ws.Cells.Item(132, 9).Value = 229240.73  # I132: 245990.12 -> 229240.73
$ws.Cells.Item(132, 11).Value = 687722.1900000001  # K132: 737970.36 -> 687722.1900000001
$ws.Cells.Item(132, 13).Value = -685192.1900000001  # M132: -735440.36 -> -685192.1900000001

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 389652.3  # H7: 375552.22 -> 389652.3
$ws.Cells.Item(7, 9).Value = 4017.5293  # I7: 4073.4119 -> 4017.5293
$ws.Cells.Item(7, 10).Value = 1118073.5  # J7: 1007066.2 -> 1118073.5
$ws.Cells.Item(7, 11).Value = 4017.5293  # K7: 4073.4119 -> 4017.5293
$ws.Cells.Item(7, 12).Value = 1118073.5  # L7: 1007066.2 -> 1118073.5
$ws.Cells.Item(7, 13).Value = -3905.5293  # M7: -3961.4119 -> -3905.5293
$ws.Cells.Item(7, 14).Value = -1118297.5  # N7: -1007290.2 -> -1118297.5
$ws.Cells.Item(22, 8).Value = 699.1667  # H22: 587.625 -> 699.1667
$ws.Cells.Item(22, 10).Value = 949.25  # J22: 976.25 -> 949.25
$ws.Cells.Item(22, 12).Value = 949.25  # L22: 976.25 -> 949.25
$ws.Cells.Item(22, 14).Value = -1539.25  # N22: -1566.25 -> -1539.25
$ws.Cells.Item(27, 8).Value = 699.1667  # H27: 587.625 -> 699.1667
$ws.Cells.Item(27, 10).Value = 949.25  # J27: 976.25 -> 949.25
$ws.Cells.Item(27, 12).Value = 949.25  # L27: 976.25 -> 949.25
$ws.Cells.Item(27, 14).Value = -1163.25  # N27: -1190.25 -> -1163.25
$ws.Cells.Item(40, 8).Value = 89650.664  # H40: 97073.45 -> 89650.664
$ws.Cells.Item(40, 10).Value = 8250  # J40: 8285.714 -> 8250
$ws.Cells.Item(40, 12).Value = 8250  # L40: 8285.714 -> 8250
$ws.Cells.Item(40, 14).Value = -8522  # N40: -8557.714 -> -8522
$ws.Cells.Item(46, 8).Value = 3604.4285  # H46: 3508.3914 -> 3604.4285
$ws.Cells.Item(46, 9).Value = 2861.9375  # I46: 2821.7222 -> 2861.9375
$ws.Cells.Item(46, 11).Value = 2861.9375  # K46: 2821.7222 -> 2861.9375
$ws.Cells.Item(46, 13).Value = -2673.9375  # M46: -2633.7222 -> -2673.9375
$ws.Cells.Item(55, 8).Value = 807.129  # H55: 831.9666999999999 -> 807.129
$ws.Cells.Item(55, 9).Value = 167.36842  # I55: 173.22223 -> 167.36842
$ws.Cells.Item(55, 11).Value = 167.36842  # K55: 173.22223 -> 167.36842
$ws.Cells.Item(55, 13).Value = 5.631580000000014  # M55: -0.2222299999999962 -> 5.631580000000014
$ws.Cells.Item(68, 8).Value = 69630.625  # H68: 73952.664 -> 69630.625
$ws.Cells.Item(68, 9).Value = 4899.75  # I68: 4933 -> 4899.75
$ws.Cells.Item(68, 11).Value = 4899.75  # K68: 4933 -> 4899.75
$ws.Cells.Item(68, 13).Value = -4150.75  # M68: -4184 -> -4150.75
$ws.Cells.Item(71, 8).Value = 69630.625  # H71: 73952.664 -> 69630.625
$ws.Cells.Item(71, 9).Value = 4899.75  # I71: 4933 -> 4899.75
$ws.Cells.Item(71, 11).Value = 24498.75  # K71: 24665 -> 24498.75
$ws.Cells.Item(71, 13).Value = -20754.75  # M71: -20921 -> -20754.75
$ws.Cells.Item(74, 8).Value = 36829.637  # H74: 36548.08 -> 36829.637
$ws.Cells.Item(74, 9).Value = 50197  # I74: 42848 -> 50197
$ws.Cells.Item(74, 10).Value = 35492.9  # J74: 35402.637 -> 35492.9
$ws.Cells.Item(74, 11).Value = 50197  # K74: 42848 -> 50197
$ws.Cells.Item(74, 12).Value = 35492.9  # L74: 35402.637 -> 35492.9
$ws.Cells.Item(74, 13).Value = -49199  # M74: -41850 -> -49199
$ws.Cells.Item(74, 14).Value = -37488.9  # N74: -37398.637 -> -37488.9
$ws.Cells.Item(77, 8).Value = 36829.637  # H77: 36548.08 -> 36829.637
$ws.Cells.Item(77, 9).Value = 50197  # I77: 42848 -> 50197
$ws.Cells.Item(77, 10).Value = 35492.9  # J77: 35402.637 -> 35492.9
$ws.Cells.Item(77, 11).Value = 150591  # K77: 128544 -> 150591
$ws.Cells.Item(77, 12).Value = 106478.7  # L77: 106207.911 -> 106478.7
$ws.Cells.Item(77, 13).Value = -145599  # M77: -123552 -> -145599
$ws.Cells.Item(77, 14).Value = -116462.7  # N77: -116191.911 -> -116462.7
$ws.Cells.Item(122, 8).Value = 1255053.4  # H122: 1181523.9 -> 1255053.4
$ws.Cells.Item(122, 9).Value = 1114316.5  # I122: 1253175.1 -> 1114316.5
$ws.Cells.Item(122, 10).Value = 1436000.8  # J122: 1117833.8 -> 1436000.8
$ws.Cells.Item(122, 11).Value = 3342949.5  # K122: 3759525.3 -> 3342949.5
$ws.Cells.Item(122, 12).Value = 4308002.4  # L122: 3353501.4 -> 4308002.4
$ws.Cells.Item(122, 13).Value = -3340499.5  # M122: -3757075.3 -> -3340499.5
$ws.Cells.Item(122, 14).Value = -4312902.4  # N122: -3358401.4 -> -4312902.4
$ws.Cells.Item(126, 8).Value = 389652.3  # H126: 375552.22 -> 389652.3
$ws.Cells.Item(126, 9).Value = 4017.5293  # I126: 4073.4119 -> 4017.5293
$ws.Cells.Item(126, 10).Value = 1118073.5  # J126: 1007066.2 -> 1118073.5
$ws.Cells.Item(126, 11).Value = 12052.5879  # K126: 12220.2357 -> 12052.5879
$ws.Cells.Item(126, 12).Value = 3354220.5  # L126: 3021198.6 -> 3354220.5
$ws.Cells.Item(126, 13).Value = -9582.5879  # M126: -9750.235700000001 -> -9582.5879
$ws.Cells.Item(126, 14).Value = -3359160.5  # N126: -3026138.6 -> -3359160.5
$ws.Cells.Item(132, 8).Value = 5283.2827  # H132: 4743.8867 -> 5283.2827
$ws.Cells.Item(132, 10).Value = 7540  # J132: 4929.1177 -> 7540
$ws.Cells.Item(132, 12).Value = 22620  # L132: 14787.3531 -> 22620
$ws.Cells.Item(132, 14).Value = -27680  # N132: -19847.3531 -> -27680
$ws.Cells.Item(136, 8).Value = 670805.3  # H136: 559248.9 -> 670805.3
$ws.Cells.Item(136, 9).Value = 955366.9399999999  # I136: 743389.1 -> 955366.9399999999
$ws.Cells.Item(136, 11).Value = 2866100.82  # K136: 2230167.3 -> 2866100.82
$ws.Cells.Item(136, 13).Value = -2863550.82  # M136: -2227617.3 -> -2863550.82

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 8).Value = 95490.44  # H40: 92935.89 -> 95490.44
$ws.Cells.Item(40, 9).Value = 95490.44  # I40: 92935.89 -> 95490.44
$ws.Cells.Item(40, 11).Value = 95490.44  # K40: 92935.89 -> 95490.44
$ws.Cells.Item(40, 13).Value = -95341.44  # M40: -92786.89 -> -95341.44
$ws.Cells.Item(100, 8).Value = 870.82355  # H100: 1021.5714 -> 870.82355
$ws.Cells.Item(100, 9).Value = 953.2  # I100: 1080.6923 -> 953.2
$ws.Cells.Item(100, 11).Value = 1906.4  # K100: 2161.3846 -> 1906.4
$ws.Cells.Item(100, 13).Value = -1365.4  # M100: -1620.3846 -> -1365.4
$ws.Cells.Item(107, 8).Value = 562.13336  # H107: 619.38464 -> 562.13336
$ws.Cells.Item(107, 9).Value = 517.3333  # I107: 610.8570999999999 -> 517.3333
$ws.Cells.Item(107, 11).Value = 1551.9999  # K107: 1832.5713 -> 1551.9999
$ws.Cells.Item(107, 13).Value = 368.0001  # M107: 87.42870000000016 -> 368.0001
$ws.Cells.Item(113, 8).Value = 773.875  # H113: 794.913 -> 773.875
$ws.Cells.Item(113, 10).Value = 1130.7778  # J113: 1235.875 -> 1130.7778
$ws.Cells.Item(113, 12).Value = 3392.3334  # L113: 3707.625 -> 3392.3334
$ws.Cells.Item(113, 14).Value = -7732.3334  # N113: -8047.625 -> -7732.3334
$ws.Cells.Item(122, 8).Value = 45459540  # H122: 34486784 -> 45459540
$ws.Cells.Item(122, 9).Value = 111113610  # I122: 66668556 -> 111113610
$ws.Cells.Item(122, 10).Value = 6726.385  # J122: 6315.5713 -> 6726.385
$ws.Cells.Item(122, 11).Value = 333340830  # K122: 200005668 -> 333340830
$ws.Cells.Item(122, 12).Value = 20179.155  # L122: 18946.7139 -> 20179.155
$ws.Cells.Item(122, 13).Value = -333338380  # M122: -200003218 -> -333338380
$ws.Cells.Item(122, 14).Value = -25079.155  # N122: -23846.7139 -> -25079.155
$ws.Cells.Item(126, 8).Value = 2319.6875  # H126: 2447.9412 -> 2319.6875
$ws.Cells.Item(126, 10).Value = 5500  # J126: 5166.6665 -> 5500
$ws.Cells.Item(126, 12).Value = 16500  # L126: 15499.9995 -> 16500
$ws.Cells.Item(126, 14).Value = -21440  # N126: -20439.9995 -> -21440
$ws.Cells.Item(132, 8).Value = 31226.166  # H132: 31164.334 -> 31226.166
$ws.Cells.Item(132, 9).Value = 2560.8518  # I132: 2497 -> 2560.8518
$ws.Cells.Item(132, 10).Value = 117222.11  # J132: 131500 -> 117222.11
$ws.Cells.Item(132, 11).Value = 7682.555399999999  # K132: 7491 -> 7682.555399999999
$ws.Cells.Item(132, 12).Value = 351666.33  # L132: 394500 -> 351666.33
$ws.Cells.Item(132, 13).Value = -5152.555399999999  # M132: -4961 -> -5152.555399999999
$ws.Cells.Item(132, 14).Value = -356726.33  # N132: -399560 -> -356726.33
$ws.Cells.Item(136, 8).Value = 367589.34  # H136: 390520.44 -> 367589.34
$ws.Cells.Item(136, 9).Value = 1042544  # I136: 1303007 -> 1042544
$ws.Cells.Item(136, 11).Value = 3127632  # K136: 3909021 -> 3127632
$ws.Cells.Item(136, 13).Value = -3125082  # M136: -3906471 -> -3125082
